$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that no longer have content (do this first so later Sets are not clobbered)
$ws.Range("A18").ClearContents()
$ws.Range("B20").ClearContents()
$ws.Range("C20").ClearContents()
$ws.Range("D20").ClearContents()
$ws.Range("E20").ClearContents()
$ws.Range("A26").ClearContents()
$ws.Range("E28").ClearContents()
$ws.Range("B29").ClearContents()
$ws.Range("D29").ClearContents()
$ws.Range("B31").ClearContents()
$ws.Range("C31").ClearContents()
$ws.Range("D31").ClearContents()
$ws.Range("E31").ClearContents()

# Update cell text content that changed (plan items re-flowed / re-worded after moving files into Documentation folder)
$ws.Range("A11").Value = "- A default list of 30 names and facial features need to generated"
$ws.Range("A13").Value = "List and dictionies data base."
$ws.Range("A14").Value = "- Create a list for names of all the player and fictional character."
$ws.Range("A15").Value = "- Creates a list of descritive words used within each facial features dictionay"
$ws.Range("A16").Value = "- Need to ensure that the descritive words aren't repeated within their lists "
$ws.Range("A17").Value = "- A default list of 20 names and facial features need to generated needed "
$ws.Range("A19").Value = "Main game feature "
$ws.Range("B19").Value = "High "
$ws.Range("C19").Value = "10 hours"
$ws.Range("D19").Value = 44391
$ws.Range("E19").Value = "Plenty of time needs to be allocated for this. "
$ws.Range("A20").Value = "-Communicate to the user the how the game works. "
$ws.Range("A21").Value = "- It will need to be able to generate randon questions"
$ws.Range("A22").Value = "- it will need to take in user responses and ensure no errors occor"
$ws.Range("A23").Value = "- needs remove the question it has asked in order not to ask it again"
$ws.Range("A24").Value = "- names to be filttered depending on the response of the user."
$ws.Range("A25").Value = "- Has a counter to figure out how many questions have been asked. "
$ws.Range("A27").Value = "Have the terminal cleared after each question "
$ws.Range("B27").Value = "low "
$ws.Range("D27").Value = 44393
$ws.Range("E27").Value = "This is a nice to have and is to be completed once all other work has been done "
$ws.Range("A28").Value = "A help file needs to be created for the users "
$ws.Range("B28").Value = "Medium "
$ws.Range("A29").Value = "-"
$ws.Range("A30").Value = "Bash script "
$ws.Range("B30").Value = "Medium "
$ws.Range("C30").Value = "2 hours "
$ws.Range("D30").Value = 44393
$ws.Range("E30").Value = "This is to be done once Python script is complete"
$ws.Range("A31").Value = "-"

# Remove the now-unused last row (content shifted up by one, dimension shrinks from E44 to E43)
$ws.Rows.Item(44).Delete()

# Update selection to match the saved cursor position
$ws.Application.Goto($ws.Range("A24"))
